# Updates the "Price" column (D) values for several cryptocurrency rows
# on Sheet1, matching the latest scrape snapshot.
# Each target cell currently stores its price as text (inline string);
# we explicitly format the cell as Text before writing so Excel keeps
# the value as a string (instead of silently converting it to a
# floating point number), then restore the cell's style to "Normal" so
# no extra formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "245.02" }
    @{ Cell = "D3"; Value = "23.85" }
    @{ Cell = "D4"; Value = "5.322" }
    @{ Cell = "D5"; Value = "0.05781" }
    @{ Cell = "D6"; Value = "6.467" }
    @{ Cell = "D7"; Value = "3.329" }
    @{ Cell = "D8"; Value = "0.8129" }
    @{ Cell = "D9"; Value = "0.8892" }
    @{ Cell = "D10"; Value = "0.1392" }
    @{ Cell = "D11"; Value = "0.07349" }
    @{ Cell = "D12"; Value = "0.03093" }
    @{ Cell = "D14"; Value = "0.09347" }
    @{ Cell = "D15"; Value = "3.865" }
    @{ Cell = "D16"; Value = "0.001539" }
    @{ Cell = "D17"; Value = "0.04716" }
    @{ Cell = "D18"; Value = "0.0006054" }
    @{ Cell = "D19"; Value = "0.005969" }
    @{ Cell = "D20"; Value = "0.001297" }
    @{ Cell = "D22"; Value = "0.00008806" }
    @{ Cell = "D23"; Value = "3.583" }
    @{ Cell = "D25"; Value = "0.3179" }
    @{ Cell = "D40"; Value = "0.03810" }
    @{ Cell = "D42"; Value = "0.002752" }
    @{ Cell = "D43"; Value = "0.003206" }
    @{ Cell = "D44"; Value = "0.007847" }
    @{ Cell = "D45"; Value = "0.00005471" }
    @{ Cell = "D47"; Value = "0.5503" }
    @{ Cell = "D48"; Value = "0.001846" }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    $range.NumberFormat = "@"
    $range.Value = $update.Value
    $range.Style = "Normal"
}
